# Commit: "fixed a mistake: I was using log10 instead of natural logarithms for the weights"
#
# Changes applied:
#  1) L28: pseudo-count parameter changes from 1 to 2 (drives the L30:S33 /
#     L40:S43 / L48:S51 dependent formulas via recalculation).
#  2) A brand-new helper cell F38 = LN(B30/V40) is introduced.
#  3) Every LOG10-based weight-matrix formula in the two score tables
#     (rows 40-43 using $V40:$V43, and rows 48-51 using $V48:$V51, each
#     covering both the B:I block and the L:S block) is rewritten from
#     LOG(...) to LN(...).
#  4) The frozen-pane top-left cell and the active selection are moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the pseudo-count parameter -----------------------------------
$ws.Range("L28").Value = 2

# --- 2) Add the new F38 helper cell --------------------------------------
$ws.Range("F38").Formula = "=LN(B30/V40)"

# --- 3) Replace LOG(...) with LN(...) in both score-weight tables --------
$dataCols  = @("B","C","D","E","F","G","H","I")
$pseudoCols = @("L","M","N","O","P","Q","R","S")

# Block 1: rows 40-43 -> data rows 30-33, pseudo-count column $V (rows 40-43)
# Block 2: rows 48-51 -> data rows 30-33, pseudo-count column $V (rows 48-51)
$blocks = @(
    @{ FormulaRows = 40..43; VCol = "V" },
    @{ FormulaRows = 48..51; VCol = "V" }
)

foreach ($block in $blocks) {
    $dataRow = 30
    foreach ($fRow in $block.FormulaRows) {
        for ($i = 0; $i -lt $dataCols.Length; $i++) {
            $dCol = $dataCols[$i]
            $pCol = $pseudoCols[$i]

            $formulaData = '=IF(' + $dCol + $dataRow + '=0, "-Inf",LN(' + $dCol + $dataRow + '/$' + $block.VCol + $fRow + '))'
            $ws.Range($dCol + $fRow).Formula = $formulaData

            $formulaPseudo = '=IF(' + $pCol + $dataRow + '=0, "-Inf",LN(' + $pCol + $dataRow + '/$' + $block.VCol + $fRow + '))'
            $ws.Range($pCol + $fRow).Formula = $formulaPseudo
        }
        $dataRow = $dataRow + 1
    }
}

# --- 4) Move the frozen pane / active selection ---------------------------
$ws.Range("A27").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A27").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("L50").Select()
